$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for A2:C21 after the edit
$data = @(
    @("寒武纪-U", "工业富联", "岩山科技"),
    @("先导智能", "国轩高科", "华胜天成"),
    @("中芯国际", "岩山科技", "工业富联"),
    @("国轩高科", "先导智能", "长城军工"),
    @("工业富联", "寒武纪-U", "领益智造"),
    @("岩山科技", "通富微电", "通富微电"),
    @("宁德时代", "中国稀土", "寒武纪"),
    @("华胜天成", "东方财富", "吉视传媒"),
    @("长城军工", "中芯国际", "三维通信"),
    @("通富微电", "长城军工", "中芯国际"),
    @("中国稀土", "华工科技", "国轩高科"),
    @("领益智造", "华胜天成", "北方稀土"),
    @("三维通信", "领益智造", "中国稀土"),
    @("吉视传媒", "三维通信", "拓维信息"),
    @("东方财富", "吉视传媒", "先导智能"),
    @("全志科技", "贵州茅台", "比亚迪"),
    @("北方稀土", "宁德时代", "利欧股份"),
    @("盛和资源", "新易盛", "天融信"),
    @("数据港", "盛和资源", "宁德时代"),
    @("润建股份", "北方稀土", "德艺文创")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
